# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, $value) {
    # Force the cell to stay a text value (matches the source data's
    # inline-string typing) even when the new text looks numeric,
    # then restore the default "Normal" style so no stray number
    # format is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '86.789.56'
Set-TextValue $ws.Range('E2') '  +7.71%  '
Set-TextValue $ws.Range('D3') '3.303.10'
Set-TextValue $ws.Range('E3') '  +4.06%  '
Set-TextValue $ws.Range('E4') '  -0.32%  '
Set-TextValue $ws.Range('D5') '217.75'
Set-TextValue $ws.Range('E5') '  +3.73%  '
Set-TextValue $ws.Range('D6') '632.75'
Set-TextValue $ws.Range('E6') '  +1.47%  '
Set-TextValue $ws.Range('D7') '0.322'
Set-TextValue $ws.Range('E7') '  +18.34%  '
Set-TextValue $ws.Range('E8') '  -0.24%  '
Set-TextValue $ws.Range('D9') '0.609'
Set-TextValue $ws.Range('E9') '  +3.68%  '
Set-TextValue $ws.Range('D10') '3.298.30'
Set-TextValue $ws.Range('E10') '  +3.80%  '
Set-TextValue $ws.Range('E11') '  +2.18%  '
Set-TextValue $ws.Range('E12') '  +5.74%  '
Set-TextValue $ws.Range('E13') '  +1.85%  '
Set-TextValue $ws.Range('D14') '3.913.13'
Set-TextValue $ws.Range('E14') '  +3.26%  '
Set-TextValue $ws.Range('E15') '  +7.94%  '
Set-TextValue $ws.Range('E16') '  +2.23%  '
Set-TextValue $ws.Range('D17') '86.636.68'
Set-TextValue $ws.Range('E17') '  +7.37%  '
Set-TextValue $ws.Range('D18') '3.300.37'
Set-TextValue $ws.Range('E18') '  +3.55%  '
Set-TextValue $ws.Range('D19') '14.41'
Set-TextValue $ws.Range('E19') '  +1.43%  '
Set-TextValue $ws.Range('D20') '3.13'
Set-TextValue $ws.Range('E20') '  +4.34%  '
Set-TextValue $ws.Range('D21') '449.58'
Set-TextValue $ws.Range('E21') '  +2.99%  '
Set-TextValue $ws.Range('D22') '9.00'
Set-TextValue $ws.Range('E22') '  -1.65%  '
Set-TextValue $ws.Range('D23') '5.29'
Set-TextValue $ws.Range('E23') '  +2.80%  '
Set-TextValue $ws.Range('D24') '7.39'
Set-TextValue $ws.Range('E24') '  +6.19%  '
Set-TextValue $ws.Range('D25') '5.31'
Set-TextValue $ws.Range('E25') '  +13.85%  '
Set-TextValue $ws.Range('D26') '12.35'
Set-TextValue $ws.Range('E26') '  +13.89%  '
Set-TextValue $ws.Range('D27') '3.495.08'
Set-TextValue $ws.Range('E27') '  +4.80%  '
Set-TextValue $ws.Range('E28') '  +2.29%  '
Set-TextValue $ws.Range('D29') '0.210'
Set-TextValue $ws.Range('E29') '  +70.73%  '
Set-TextValue $ws.Range('B30') 'PEPE'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D30') '0.0000127'
Set-TextValue $ws.Range('E30') '  +5.15%  '
Set-TextValue $ws.Range('B31') 'Dai'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D31') '0.999'
Set-TextValue $ws.Range('E31') '  -0.02%  '
Set-TextValue $ws.Range('D32') '9.19'
Set-TextValue $ws.Range('E32') '  +3.07%  '
Set-TextValue $ws.Range('D33') '588.68'
Set-TextValue $ws.Range('E33') '  +5.68%  '
Set-TextValue $ws.Range('E34') '  +0.54%  '
Set-TextValue $ws.Range('D35') '1.52'
Set-TextValue $ws.Range('E35') '  +3.67%  '
Set-TextValue $ws.Range('E36') '  +2.42%  '
Set-TextValue $ws.Range('E37') '  +0.04%  '
Set-TextValue $ws.Range('D38') '23.34'
Set-TextValue $ws.Range('E38') '  +1.79%  '
Set-TextValue $ws.Range('D39') '6.52'
Set-TextValue $ws.Range('E39') '  +15.48%  '
Set-TextValue $ws.Range('D40') '0.998'
Set-TextValue $ws.Range('E40') '  -0.43%  '
Set-TextValue $ws.Range('D41') '0.414'
Set-TextValue $ws.Range('E41') '  +2.85%  '
Set-TextValue $ws.Range('E42') '  +2.88%  '
Set-TextValue $ws.Range('E43') '  +13.22%  '
Set-TextValue $ws.Range('D44') '3.02'
Set-TextValue $ws.Range('E44') '  +12.78%  '
Set-TextValue $ws.Range('D45') '158.43'
Set-TextValue $ws.Range('E45') '  -3.68%  '
Set-TextValue $ws.Range('D47') '187.75'
Set-TextValue $ws.Range('E47') '  -1.11%  '
Set-TextValue $ws.Range('D48') '46.51'
Set-TextValue $ws.Range('E48') '  +8.55%  '
Set-TextValue $ws.Range('D49') '1.34'
Set-TextValue $ws.Range('E49') '  +3.96%  '
Set-TextValue $ws.Range('E50') '  -0.15%  '
Set-TextValue $ws.Range('D51') '26.17'
Set-TextValue $ws.Range('E51') '  +6.65%  '
